$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "runs/balls/fours/sixes" figures for the Marcus Stoinis innings
# rows (activity re-synced from the Excel source form). Only the cells
# whose numbers actually changed are touched; values are kept as text
# (matching the sheet's existing "number stored as text" convention).
$updates = [ordered]@{
    "C2"  = "10"
    "D2"  = "5"
    "E2"  = "0"
    "F2"  = "1"
    "C3"  = "2"
    "D3"  = "3"
    "E3"  = "0"
    "C4"  = "5"
    "D4"  = "6"
    "E4"  = "1"
    "F4"  = "0"
    "C6"  = "65"
    "D6"  = "46"
    "E6"  = "6"
    "F6"  = "3"
    "C8"  = "18"
    "D8"  = "19"
    "C9"  = "24"
    "D9"  = "14"
    "E9"  = "1"
    "F9"  = "2"
    "C10" = "13"
    "D10" = "8"
    "E10" = "2"
    "C11" = "9"
    "D11" = "10"
    "E11" = "0"
    "C12" = "0"
    "D12" = "1"
    "E12" = "0"
    "F12" = "0"
    "C13" = "11"
    "D13" = "9"
    "E13" = "1"
    "F13" = "0"
    "C14" = "1"
    "D14" = "3"
    "E14" = "0"
    "C15" = "5"
    "D15" = "3"
    "E15" = "1"
    "F15" = "0"
    "C16" = "53"
    "D16" = "21"
    "E16" = "7"
    "F16" = "3"
    "C17" = "53"
    "D17" = "26"
    "F17" = "2"
    "C18" = "39"
    "D18" = "30"
    "F18" = "4"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
